$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.648.90"
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").Value = "1.755.62"
$ws.Range("E3").Value = "  -0.54%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'324.26"
$ws.Range("E5").Value = "  +1.21%  "
$ws.Range("D6").Value = "'0.9998"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4599"
$ws.Range("E7").Value = "  +7.87%  "
$ws.Range("D8").Value = "'0.3598"
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("D9").Value = "'0.07535"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'42.28"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("D11").Value = "'1.100"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("E12").Value = "  -0.10%  "
$ws.Range("D13").Value = "'20.86"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "'6.024"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "'7.115"
$ws.Range("E15").Value = "  -2.79%  "
$ws.Range("D16").Value = "1.749.44"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "'92.78"
$ws.Range("E17").Value = "  +1.76%  "
$ws.Range("E18").Value = "  +0.91%  "
$ws.Range("D19").Value = "'0.06419"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'16.80"
$ws.Range("E21").Value = "  -1.40%  "
$ws.Range("D22").Value = "'5.836"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").Value = "27.694.24"
$ws.Range("E23").Value = "  -1.19%  "
$ws.Range("D24").Value = "'11.25"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.109"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'163.84"
$ws.Range("E26").Value = "  +4.18%  "
$ws.Range("D27").Value = "'20.45"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "1.955.68"
$ws.Range("E28").Value = "  -1.24%  "
$ws.Range("D29").Value = "'2.095"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "'127.26"
$ws.Range("E30").Value = "  +2.06%  "
$ws.Range("D31").Value = "'1.077"
$ws.Range("E31").Value = "  -6.75%  "
$ws.Range("E32").Value = "  +3.97%  "
$ws.Range("D33").Value = "'3.668"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").Value = "'5.545"
$ws.Range("E34").Value = "  -1.44%  "
$ws.Range("E35").Value = "  -4.49%  "
$ws.Range("D36").Value = "'0.02300"
$ws.Range("E36").Value = "  -0.23%  "
$ws.Range("D37").Value = "'0.2102"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'0.06044"
$ws.Range("E38").Value = "  +0.21%  "
$ws.Range("D39").Value = "'0.6383"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").Value = "'4.972"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("D41").Value = "'1.202"
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("E42").Value = "  -1.14%  "
$ws.Range("D43").Value = "'7.828"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("D44").Value = "'13.32"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "'0.5920"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'3.710"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "'123.23"
$ws.Range("E47").Value = "  +0.70%  "
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "'1.149"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("D50").Value = "'0.06865"
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'72.32"
$ws.Range("E51").Value = "  -2.21%  "
